$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Global text substitutions shared by every sheet (file renames + status).
# ---------------------------------------------------------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("ac2f5080-7f2a-49cf-be49-8ed8dacc307e.md", "5a57b781-858c-4266-b122-ad0635dcfa74.md")
    $ws.Cells.Replace("afb46ee5-6896-4257-a56d-04be8f8c5f92.md", "ffffdb97a82b-372e-4208-ad5f-b8e8181f25c3.md")
    $ws.Cells.Replace("Handed back: in sync with en-US", "Ready for handoff")
    $ws.Cells.Replace("2016-08-23 00:57:52", "2016-08-23 00:59:04")
}

# ---------------------------------------------------------------------------
# 2. Overview sheet: fix hyperlink display text (targets are unchanged).
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/ac2f5080-7f2a-49cf-be49-8ed8dacc307e.md", "", "", "e2e\5a57b781-858c-4266-b122-ad0635dcfa74.md")
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/afb46ee5-6896-4257-a56d-04be8f8c5f92.md", "", "", "e2e\ffffdb97a82b-372e-4208-ad5f-b8e8181f25c3.md")

# ---------------------------------------------------------------------------
# 3. zh-cn sheet (sheet2): xliff handoff is now complete -> clear the
#    "Latest Target File" / "Latest Handback File" columns (I/J), flip
#    Content Duplicate (F3) to True, refresh the handoff file/timestamp
#    columns (G/H) and blank the handback timestamp (K) to the epoch value.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("F3").Value = "'True"
$wsZh.Range("F3").Style = "Normal"
$wsZh.Range("G2").Value = "5a57b781-858c-4266-b122-ad0635dcfa74.92d3c8e7f0eaf37956d44926cd72a9d6ac625084.zh-cn.xlf"
$wsZh.Range("G3").Value = "5a57b781-858c-4266-b122-ad0635dcfa74.92d3c8e7f0eaf37956d44926cd72a9d6ac625084.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-08-23 00:58:56"
$wsZh.Range("H3").Value = "2016-08-23 00:58:56"
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("I2").Value = ""
$wsZh.Range("J2").Value = ""
$wsZh.Range("I3").Value = ""
$wsZh.Range("J3").Value = ""
$wsZh.Range("I2:J3").ClearFormats()

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/ac2f5080-7f2a-49cf-be49-8ed8dacc307e.md", "", "", "5a57b781-858c-4266-b122-ad0635dcfa74.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/afb46ee5-6896-4257-a56d-04be8f8c5f92.md", "", "", "ffffdb97a82b-372e-4208-ad5f-b8e8181f25c3.md")

# ---------------------------------------------------------------------------
# 4. de-de sheet (sheet3): same shape of change as zh-cn above.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("F3").Value = "'True"
$wsDe.Range("F3").Style = "Normal"
$wsDe.Range("G2").Value = "5a57b781-858c-4266-b122-ad0635dcfa74.92d3c8e7f0eaf37956d44926cd72a9d6ac625084.de-de.xlf"
$wsDe.Range("G3").Value = "5a57b781-858c-4266-b122-ad0635dcfa74.92d3c8e7f0eaf37956d44926cd72a9d6ac625084.de-de.xlf"
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("I2").Value = ""
$wsDe.Range("J2").Value = ""
$wsDe.Range("I3").Value = ""
$wsDe.Range("J3").Value = ""
$wsDe.Range("I2:J3").ClearFormats()

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/ac2f5080-7f2a-49cf-be49-8ed8dacc307e.md", "", "", "5a57b781-858c-4266-b122-ad0635dcfa74.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3da8bc095bc12110e89b16dcad1b50726d1bdebd/e2e/afb46ee5-6896-4257-a56d-04be8f8c5f92.md", "", "", "ffffdb97a82b-372e-4208-ad5f-b8e8181f25c3.md")

# ---------------------------------------------------------------------------
# 5. Column width tweaks (best-effort; engine quantizes to 1/6-character
#    steps so we pick the closest representable width to the target).
# ---------------------------------------------------------------------------
$wsOverview.Columns("E:F").ColumnWidth = 16.38
$wsZh.Columns("C:C").ColumnWidth = 16.38
$wsZh.Columns("I:I").ColumnWidth = 17.82
$wsZh.Columns("J:J").ColumnWidth = 20.87
$wsDe.Columns("C:C").ColumnWidth = 16.38
$wsDe.Columns("I:I").ColumnWidth = 17.82
$wsDe.Columns("J:J").ColumnWidth = 20.87

Write-Host "edit.ps1 completed"
